{"js": "// Update the three-digit \u00f7 one-digit division worksheet answers.\n//\n// The document body contains a single table of \"XXX\u00f7Y=ZZ, R\" answer\n// strings; every populated cell's text is unique across the whole\n// document, so an exact-text search/replace on `context.document.body`\n// unambiguously targets the correct run without needing to walk table\n// rows/cells explicitly.\nconst replacements = [\n  [\"511\u00f73=170, 1\", \"469\u00f72=234, 1\"],\n  [\"396\u00f74=99, 0\", \"865\u00f79=96, 1\"],\n  [\"746\u00f78=93, 2\", \"813\u00f76=135, 3\"],\n  [\"504\u00f77=72, 0\", \"282\u00f73=94, 0\"],\n  [\"307\u00f75=61, 2\", \"980\u00f73=326, 2\"],\n  [\"613\u00f76=102, 1\", \"887\u00f73=295, 2\"],\n  [\"132\u00f75=26, 2\", \"804\u00f73=268, 0\"],\n  [\"496\u00f77=70, 6\", \"160\u00f79=17, 7\"],\n  [\"549\u00f73=183, 0\", \"669\u00f76=111, 3\"],\n  [\"486\u00f78=60, 6\", \"290\u00f76=48, 2\"],\n  [\"452\u00f79=50, 2\", \"716\u00f79=79, 5\"],\n  [\"744\u00f74=186, 0\", \"694\u00f72=347, 0\"],\n  [\"707\u00f72=353, 1\", \"958\u00f77=136, 6\"],\n  [\"422\u00f75=84, 2\", \"587\u00f73=195, 2\"],\n  [\"482\u00f78=60, 2\", \"535\u00f74=133, 3\"],\n  [\"596\u00f79=66, 2\", \"298\u00f74=74, 2\"],\n  [\"617\u00f79=68, 5\", \"126\u00f78=15, 6\"],\n  [\"249\u00f74=62, 1\", \"519\u00f76=86, 3\"],\n  [\"199\u00f74=49, 3\", \"262\u00f79=29, 1\"],\n  [\"486\u00f72=243, 0\", \"209\u00f73=69, 2\"],\n  [\"105\u00f77=15, 0\", \"606\u00f78=75, 6\"],\n  [\"619\u00f73=206, 1\", \"155\u00f77=22, 1\"],\n  [\"329\u00f74=82, 1\", \"298\u00f79=33, 1\"],\n  [\"296\u00f73=98, 2\", \"440\u00f77=62, 6\"],\n  [\"214\u00f79=23, 7\", \"931\u00f76=155, 1\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items.forEach((range) => {\n    range.insertText(newText, Word.InsertLocation.replace);\n  });\n}\n\nawait context.sync();\n", "ps1": "# Update the three-digit \u00f7 one-digit division worksheet answers.\n#\n# The document body contains a single table of \"XXX\u00f7Y=ZZ, R\" answer\n# strings; every populated cell's text is unique across the whole\n# document, so exact-text Find/Replace unambiguously targets the\n# correct run without needing to walk table rows/cells explicitly.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"511\u00f73=170, 1\"; New = \"469\u00f72=234, 1\" },\n    @{ Old = \"396\u00f74=99, 0\"; New = \"865\u00f79=96, 1\" },\n    @{ Old = \"746\u00f78=93, 2\"; New = \"813\u00f76=135, 3\" },\n    @{ Old = \"504\u00f77=72, 0\"; New = \"282\u00f73=94, 0\" },\n    @{ Old = \"307\u00f75=61, 2\"; New = \"980\u00f73=326, 2\" },\n    @{ Old = \"613\u00f76=102, 1\"; New = \"887\u00f73=295, 2\" },\n    @{ Old = \"132\u00f75=26, 2\"; New = \"804\u00f73=268, 0\" },\n    @{ Old = \"496\u00f77=70, 6\"; New = \"160\u00f79=17, 7\" },\n    @{ Old = \"549\u00f73=183, 0\"; New = \"669\u00f76=111, 3\" },\n    @{ Old = \"486\u00f78=60, 6\"; New = \"290\u00f76=48, 2\" },\n    @{ Old = \"452\u00f79=50, 2\"; New = \"716\u00f79=79, 5\" },\n    @{ Old = \"744\u00f74=186, 0\"; New = \"694\u00f72=347, 0\" },\n    @{ Old = \"707\u00f72=353, 1\"; New = \"958\u00f77=136, 6\" },\n    @{ Old = \"422\u00f75=84, 2\"; New = \"587\u00f73=195, 2\" },\n    @{ Old = \"482\u00f78=60, 2\"; New = \"535\u00f74=133, 3\" },\n    @{ Old = \"596\u00f79=66, 2\"; New = \"298\u00f74=74, 2\" },\n    @{ Old = \"617\u00f79=68, 5\"; New = \"126\u00f78=15, 6\" },\n    @{ Old = \"249\u00f74=62, 1\"; New = \"519\u00f76=86, 3\" },\n    @{ Old = \"199\u00f74=49, 3\"; New = \"262\u00f79=29, 1\" },\n    @{ Old = \"486\u00f72=243, 0\"; New = \"209\u00f73=69, 2\" },\n    @{ Old = \"105\u00f77=15, 0\"; New = \"606\u00f78=75, 6\" },\n    @{ Old = \"619\u00f73=206, 1\"; New = \"155\u00f77=22, 1\" },\n    @{ Old = \"329\u00f74=82, 1\"; New = \"298\u00f79=33, 1\" },\n    @{ Old = \"296\u00f73=98, 2\"; New = \"440\u00f77=62, 6\" },\n    @{ Old = \"214\u00f79=23, 7\"; New = \"931\u00f76=155, 1\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $found = $range.Find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n    if (-not $found) {\n        throw \"Text not found: $($pair.Old)\"\n    }\n}\n"}
